# Update countries & provincias Spain
# - Swap rank order of Barein/Noruega (Barein now ranks above Noruega)
# - Swap rank order of Angola/Polinesia Francesa (Angola now ranks above Polinesia Francesa)
# - Update the "Datos actualizados" timestamp string
# - Refresh several countries' statistics (rows 4,8,14,49,54,55,71,91,177,178)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder country labels so the ranking (sorted by Casos totales desc) stays correct ---
# Row 54 was Noruega, row 55 was Barein; Barein's updated totals now exceed Noruega's,
# so Barein moves up to row 54 and Noruega drops to row 55 (text swap only; the numeric
# data for each row is updated separately below to reflect each country's real figures).
$ws.Cells.Item(54,1).Value = "Barein"
$ws.Cells.Item(55,1).Value = "Noruega"

# Row 177 was Polinesia Francesa, row 178 was Angola; Angola's updated totals now
# exceed Polinesia Francesa's, so Angola moves up to row 177.
$ws.Cells.Item(177,1).Value = "Angola"
$ws.Cells.Item(178,1).Value = "Polinesia Francesa"

# --- Update the timestamp footer string ---
$ws.Cells.Item(1,1).Value = "Datos actualizados a 22 de Mayo de 2020 a las 15:35"

# --- Refresh statistics for updated rows ---
# Row 4: Estados Unidos
$ws.Cells.Item(4,2).Value = 1622191
$ws.Cells.Item(4,3).Value = 1289
$ws.Cells.Item(4,4).Value = 382936
$ws.Cells.Item(4,5).Value = 1142870
$ws.Cells.Item(4,7).Value = 31
$ws.Cells.Item(4,8).Value = 96385

# Row 8
$ws.Cells.Item(8,2).Value = 254195
$ws.Cells.Item(8,3).Value = 3287
$ws.Cells.Item(8,7).Value = 351
$ws.Cells.Item(8,8).Value = 36393

# Row 14
$ws.Cells.Item(14,2).Value = 120532
$ws.Cells.Item(14,3).Value = 2306
$ws.Cells.Item(14,4).Value = 49872
$ws.Cells.Item(14,5).Value = 67055
$ws.Cells.Item(14,7).Value = 21
$ws.Cells.Item(14,8).Value = 3605

# Row 49
$ws.Cells.Item(49,2).Value = 11024
$ws.Cells.Item(49,3).Value = 105
$ws.Cells.Item(49,5).Value = 5417

# Row 54 (now Barein, new figures)
$ws.Cells.Item(54,2).Value = 8338
$ws.Cells.Item(54,3).Value = 164
$ws.Cells.Item(54,4).Value = 4004
$ws.Cells.Item(54,5).Value = 4322
$ws.Cells.Item(54,8).Value = 12

# Row 55 (now Noruega, keeps its previous figures)
$ws.Cells.Item(55,2).Value = 8309
$ws.Cells.Item(55,4).Value = 32
$ws.Cells.Item(55,5).Value = 8042
$ws.Cells.Item(55,8).Value = 235

# Row 71
$ws.Cells.Item(71,2).Value = 3855
$ws.Cells.Item(71,3).Value = 106
$ws.Cells.Item(71,4).Value = 2399
$ws.Cells.Item(71,5).Value = 1410
$ws.Cells.Item(71,7).Value = 2
$ws.Cells.Item(71,8).Value = 46

# Row 91
$ws.Cells.Item(91,4).Value = 1791
$ws.Cells.Item(91,5).Value = 2

# Row 177 (now Angola, new figures)
$ws.Cells.Item(177,3).Value = 2
$ws.Cells.Item(177,4).Value = 17
$ws.Cells.Item(177,5).Value = 40
$ws.Cells.Item(177,8).Value = 3

# Row 178 (now Polinesia Francesa, keeps its previous figures)
$ws.Cells.Item(178,2).Value = 60
$ws.Cells.Item(178,4).Value = 60
$ws.Cells.Item(178,5).Value = 0
$ws.Cells.Item(178,8).Value = 0
